$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the sub-bullet paragraphs that were removed from the Table of
#    Contents list (the whole block collapses down to just the five bold
#    section headers: Introduction, Data visualization, Initial Hypotheses,
#    Model Fitting, Conclusion).
# ---------------------------------------------------------------------------
$targets = @(
    "Vision",
    "Mission Statementa",
    "Organizational Values",
    "Business Model",
    "Business Growth Opportunities",
    "Data visualization",
    "Industry structure ",
    "Market trends",
    "Market size and Growth Projections",
    "Target Market / Customer Profile",
    "Competitive Analysis",
    "Overall Market Strategy",
    "Positioning",
    "Product Strategy",
    "Pricing Strategy",
    "Distribution Strategy",
    "Advertising and Promotion",
    "Facilities and Premises",
    "Equipment and Production",
    "Value Chain",
    "Production Processes",
    "Organizational Structure",
    "Operating Procedures",
    "Human Resources",
    "Skills Development"
)

$toDelete = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($targets -contains $t) {
        [void]$toDelete.Add($p.Range)
    }
}

for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $toDelete[$i].Delete()
}

# ---------------------------------------------------------------------------
# 2. Turn the old "Business Description" sub-bullet into the new bold
#    "Data visualization" section header (ilvl 0, bold, purple accent1,
#    28-half-point / 14pt text).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Business Description") {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:pPr>' +
                 '<w:pStyle w:val="ListParagraph"/>' +
                 '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' +
                 '<w:tabs><w:tab w:val="left" w:pos="10260"/></w:tabs>' +
                 '<w:spacing w:line="276" w:lineRule="auto"/>' +
                 '<w:rPr><w:b/><w:color w:val="743594" w:themeColor="accent1"/><w:sz w:val="28"/></w:rPr>' +
               '</w:pPr>' +
               '<w:r>' +
                 '<w:rPr><w:b/><w:color w:val="743594" w:themeColor="accent1"/><w:sz w:val="28"/></w:rPr>' +
                 '<w:t>Data visualization</w:t>' +
               '</w:r>' +
               '</w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the last intro paragraph to
#    the very start of the page-break paragraph that follows the TOC list.
#    Adding a bookmark with the same name relocates it (removes the old one).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -eq 2 -and [int][char]$t[0] -eq 12) {
        $startPos = $p.Range.Start
        $rng = $d.Range($startPos, $startPos)
        $d.Bookmarks.Add("_GoBack", $rng)
        break
    }
}
